$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: J1 "I would like to close my server account" ---
$ws.Range("J1").Value = "I would like to close my server account"

# --- Row 4: Wanda Watermelon wants to close her server account ---
$ws.Range("A3:D3").Copy($ws.Range("A4:D4"))
$ws.Range("H3").Copy($ws.Range("H4"))
$ws.Rows.Item(4).RowHeight = 13.8

$ws.Range("A4").Value = 44180.7005092593
$ws.Range("B4").Value = "wwatermelon@example.com"
$ws.Range("C4").Value = "Wanda"
$ws.Range("D4").Value = "Watermelon"
$ws.Range("H4").Value = "Yes"
$ws.Range("J4").Value = "Yes"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:wwatermelon@example.com", [Type]::Missing, [Type]::Missing, "wwatermelon@example.com")
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("B4").Value = "wwatermelon@example.com"

# --- Row 5: Jack Jackfruit wants to close his server account too ---
$ws.Range("A3:D3").Copy($ws.Range("A5:D5"))
$ws.Range("H3").Copy($ws.Range("H5"))
$ws.Rows.Item(5).RowHeight = 13.8

$ws.Range("A5").Value = 44211.7005092593
$ws.Range("B5").Value = "jjackfruit@example.com"
$ws.Range("C5").Value = "Jack"
$ws.Range("D5").Value = "Jackfruit"
$ws.Range("H5").Value = "Yes"
$ws.Range("J5").Value = "Yes"

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:jjackfruit@example.com", [Type]::Missing, [Type]::Missing, "jjackfruit@example.com")
$ws.Range("B3").Copy($ws.Range("B5"))
$ws.Range("B5").Value = "jjackfruit@example.com"

Write-Output "done"
